$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (A1:D1) to the new short machine-friendly names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/el/la/los/las/y) inside
#    state/municipality names, e.g. "Amatenango de la Frontera" ->
#    "Amatenango De La Frontera". Each pair below is a unique exact-text
#    replacement so Replace() only touches the intended cells.
[void]$ws.Cells.Replace('Amatenango de la Frontera', 'Amatenango De La Frontera')
[void]$ws.Cells.Replace('Bejucal de Ocampo', 'Bejucal De Ocampo')
[void]$ws.Cells.Replace('Chiapa de Corzo', 'Chiapa De Corzo')
[void]$ws.Cells.Replace('Comitán de Domínguez', 'Comitán De Domínguez')
[void]$ws.Cells.Replace('Mazapa de Madero', 'Mazapa De Madero')
[void]$ws.Cells.Replace('Montecristo de Guerrero', 'Montecristo De Guerrero')
[void]$ws.Cells.Replace('Ocozocoautla de Espinosa', 'Ocozocoautla De Espinosa')
[void]$ws.Cells.Replace('Salto de Agua', 'Salto De Agua')
[void]$ws.Cells.Replace('San Cristóbal de las Casas', 'San Cristóbal De Las Casas')
[void]$ws.Cells.Replace('Hidalgo del Parral', 'Hidalgo Del Parral')
[void]$ws.Cells.Replace('Ciudad de México', 'Ciudad De México')
[void]$ws.Cells.Replace('Cuajimalpa de Morelos', 'Cuajimalpa De Morelos')
[void]$ws.Cells.Replace('San Juan del Río', 'San Juan Del Río')
[void]$ws.Cells.Replace('Estado de México', 'Estado De México')
[void]$ws.Cells.Replace('Acambay de Ruíz Castañeda', 'Acambay De Ruíz Castañeda')
[void]$ws.Cells.Replace('Almoloya de Alquisiras', 'Almoloya De Alquisiras')
[void]$ws.Cells.Replace('Almoloya de Juárez', 'Almoloya De Juárez')
[void]$ws.Cells.Replace('Atizapán de Zaragoza', 'Atizapán De Zaragoza')
[void]$ws.Cells.Replace('Chapa de Mota', 'Chapa De Mota')
[void]$ws.Cells.Replace('Coacalco de Berriozábal', 'Coacalco De Berriozábal')
[void]$ws.Cells.Replace('Ecatepec de Morelos', 'Ecatepec De Morelos')
[void]$ws.Cells.Replace('Ixtapan de la Sal', 'Ixtapan De La Sal')
[void]$ws.Cells.Replace('Naucalpan de Juárez', 'Naucalpan De Juárez')
[void]$ws.Cells.Replace('San Felipe del Progreso', 'San Felipe Del Progreso')
[void]$ws.Cells.Replace('San Martín de las Pirámides', 'San Martín De Las Pirámides')
[void]$ws.Cells.Replace('Tenango del Aire', 'Tenango Del Aire')
[void]$ws.Cells.Replace('Tenango del Valle', 'Tenango Del Valle')
[void]$ws.Cells.Replace('Tlalnepantla de Baz', 'Tlalnepantla De Baz')
[void]$ws.Cells.Replace('Valle de Chalco Solidaridad', 'Valle De Chalco Solidaridad')
[void]$ws.Cells.Replace('Villa de Allende', 'Villa De Allende')
[void]$ws.Cells.Replace('Villa del Carbón', 'Villa Del Carbón')
[void]$ws.Cells.Replace('San Miguel de Allende', 'San Miguel De Allende')
[void]$ws.Cells.Replace('Dolores Hidalgo Cuna de la Independencia Nacional', 'Dolores Hidalgo Cuna De La Independencia Nacional')
[void]$ws.Cells.Replace('Purísima del Rincón', 'Purísima Del Rincón')
[void]$ws.Cells.Replace('San Diego de la Unión', 'San Diego De La Unión')
[void]$ws.Cells.Replace('San Luis de la Paz', 'San Luis De La Paz')
[void]$ws.Cells.Replace('Santa Cruz de Juventino Rosas', 'Santa Cruz De Juventino Rosas')
[void]$ws.Cells.Replace('Silao de la Victoria', 'Silao De La Victoria')
[void]$ws.Cells.Replace('Valle de Santiago', 'Valle De Santiago')
[void]$ws.Cells.Replace('Acapulco de Juárez', 'Acapulco De Juárez')
[void]$ws.Cells.Replace('Ajuchitlán del Progreso', 'Ajuchitlán Del Progreso')
[void]$ws.Cells.Replace('Alcozauca de Guerrero', 'Alcozauca De Guerrero')
[void]$ws.Cells.Replace('Atenango del Río', 'Atenango Del Río')
[void]$ws.Cells.Replace('Atlamajalcingo del Monte', 'Atlamajalcingo Del Monte')
[void]$ws.Cells.Replace('Atoyac de Álvarez', 'Atoyac De Álvarez')
[void]$ws.Cells.Replace('Ayutla de los Libres', 'Ayutla De Los Libres')
[void]$ws.Cells.Replace('Chilapa de Álvarez', 'Chilapa De Álvarez')
[void]$ws.Cells.Replace('Chilpancingo de los Bravo', 'Chilpancingo De Los Bravo')
[void]$ws.Cells.Replace('Coahuayutla de José María Izazaga', 'Coahuayutla De José María Izazaga')
[void]$ws.Cells.Replace('Coyuca de Benítez', 'Coyuca De Benítez')
[void]$ws.Cells.Replace('Coyuca de Catalán', 'Coyuca De Catalán')
[void]$ws.Cells.Replace('Cuetzala del Progreso', 'Cuetzala Del Progreso')
[void]$ws.Cells.Replace('Cutzamala de Pinzón', 'Cutzamala De Pinzón')
[void]$ws.Cells.Replace('Huitzuco de los Figueroa', 'Huitzuco De Los Figueroa')
[void]$ws.Cells.Replace('Iguala de la Independencia', 'Iguala De La Independencia')
[void]$ws.Cells.Replace('Ixcateopan de Cuauhtémoc', 'Ixcateopan De Cuauhtémoc')
[void]$ws.Cells.Replace('Zihuatanejo de Azueta', 'Zihuatanejo De Azueta')
[void]$ws.Cells.Replace('Mártir de Cuilapan', 'Mártir De Cuilapan')
[void]$ws.Cells.Replace('Taxco de Alarcón', 'Taxco De Alarcón')
[void]$ws.Cells.Replace('Técpan de Galeana', 'Técpan De Galeana')
[void]$ws.Cells.Replace('Tepecoacuilco de Trujano', 'Tepecoacuilco De Trujano')
[void]$ws.Cells.Replace('Tixtla de Guerrero', 'Tixtla De Guerrero')
[void]$ws.Cells.Replace('Tlalixtaquilla de Maldonado', 'Tlalixtaquilla De Maldonado')
[void]$ws.Cells.Replace('Tlapa de Comonfort', 'Tlapa De Comonfort')
[void]$ws.Cells.Replace('Agua Blanca de Iturbide', 'Agua Blanca De Iturbide')
[void]$ws.Cells.Replace('Atotonilco el Grande', 'Atotonilco El Grande')
[void]$ws.Cells.Replace('Cuautepec de Hinojosa', 'Cuautepec De Hinojosa')
[void]$ws.Cells.Replace('Huasca de Ocampo', 'Huasca De Ocampo')
[void]$ws.Cells.Replace('Huejutla de Reyes', 'Huejutla De Reyes')
[void]$ws.Cells.Replace('Jacala de Ledezma', 'Jacala De Ledezma')
[void]$ws.Cells.Replace('Mineral del Chico', 'Mineral Del Chico')
[void]$ws.Cells.Replace('Mixquiahuala de Juárez', 'Mixquiahuala De Juárez')
[void]$ws.Cells.Replace('Omitlán de Juárez', 'Omitlán De Juárez')
[void]$ws.Cells.Replace('Pachuca de Soto', 'Pachuca De Soto')
[void]$ws.Cells.Replace('Progreso de Obregón', 'Progreso De Obregón')
[void]$ws.Cells.Replace('Santiago de Anaya', 'Santiago De Anaya')
[void]$ws.Cells.Replace('Tepehuacán de Guerrero', 'Tepehuacán De Guerrero')
[void]$ws.Cells.Replace('Tepeji del Río de Ocampo', 'Tepeji Del Río De Ocampo')
[void]$ws.Cells.Replace('Tula de Allende', 'Tula De Allende')
[void]$ws.Cells.Replace('Tulancingo de Bravo', 'Tulancingo De Bravo')
[void]$ws.Cells.Replace('Villa de Tezontepec', 'Villa De Tezontepec')
[void]$ws.Cells.Replace('Zacualtipán de Ángeles', 'Zacualtipán De Ángeles')
[void]$ws.Cells.Replace('Zapotlán de Juárez', 'Zapotlán De Juárez')
[void]$ws.Cells.Replace('Ahualulco de Mercado', 'Ahualulco De Mercado')
[void]$ws.Cells.Replace('Autlán de Navarro', 'Autlán De Navarro')
[void]$ws.Cells.Replace('Encarnación de Díaz', 'Encarnación De Díaz')
[void]$ws.Cells.Replace('Lagos de Moreno', 'Lagos De Moreno')
[void]$ws.Cells.Replace('Tamazula de Gordiano', 'Tamazula De Gordiano')
[void]$ws.Cells.Replace('Tepatitlán de Morelos', 'Tepatitlán De Morelos')
[void]$ws.Cells.Replace('Unión de San Antonio', 'Unión De San Antonio')
[void]$ws.Cells.Replace('Unión de Tula', 'Unión De Tula')
[void]$ws.Cells.Replace('Zacoalco de Torres', 'Zacoalco De Torres')
[void]$ws.Cells.Replace('Zapotlán el Grande', 'Zapotlán El Grande')
[void]$ws.Cells.Replace('Coalcomán de Vázquez Pallares', 'Coalcomán De Vázquez Pallares')
[void]$ws.Cells.Replace('Jonacatepec de Leandro Valle', 'Jonacatepec De Leandro Valle')
[void]$ws.Cells.Replace('Puente de Ixtla', 'Puente De Ixtla')
[void]$ws.Cells.Replace('Tetela del Volcán', 'Tetela Del Volcán')
[void]$ws.Cells.Replace('Tlaltizapán de Zapata', 'Tlaltizapán De Zapata')
[void]$ws.Cells.Replace('Zacualpan de Amilpas', 'Zacualpan De Amilpas')
[void]$ws.Cells.Replace('San Nicolás de los Garza', 'San Nicolás De Los Garza')
[void]$ws.Cells.Replace('Acatlán de Pérez Figueroa', 'Acatlán De Pérez Figueroa')
[void]$ws.Cells.Replace('Chalcatongo de Hidalgo', 'Chalcatongo De Hidalgo')
[void]$ws.Cells.Replace('Ciénega de Zimatlán', 'Ciénega De Zimatlán')
[void]$ws.Cells.Replace('Coicoyán de las Flores', 'Coicoyán De Las Flores')
[void]$ws.Cells.Replace('Constancia del Rosario', 'Constancia Del Rosario')
[void]$ws.Cells.Replace('Cuilápam de Guerrero', 'Cuilápam De Guerrero')
[void]$ws.Cells.Replace('El Barrio de la Soledad', 'El Barrio De La Soledad')
[void]$ws.Cells.Replace('Eloxochitlán de Flores Magón', 'Eloxochitlán De Flores Magón')
[void]$ws.Cells.Replace('Fresnillo de Trujano', 'Fresnillo De Trujano')
[void]$ws.Cells.Replace('Heroica Ciudad de Ejutla de Crespo', 'Heroica Ciudad De Ejutla De Crespo')
[void]$ws.Cells.Replace('Heroica Ciudad de Huajuapan de León', 'Heroica Ciudad De Huajuapan De León')
[void]$ws.Cells.Replace('Heroica Ciudad de Tlaxiaco', 'Heroica Ciudad De Tlaxiaco')
[void]$ws.Cells.Replace('Ixtlán de Juárez', 'Ixtlán De Juárez')
[void]$ws.Cells.Replace('Heroica Ciudad de Juchitán de Zaragoza', 'Heroica Ciudad De Juchitán De Zaragoza')
[void]$ws.Cells.Replace('Mariscala de Juárez', 'Mariscala De Juárez')
[void]$ws.Cells.Replace('Mártires de Tacubaya', 'Mártires De Tacubaya')
[void]$ws.Cells.Replace('Mazatlán Villa de Flores', 'Mazatlán Villa De Flores')
[void]$ws.Cells.Replace('Miahuatlán de Porfirio Díaz', 'Miahuatlán De Porfirio Díaz')
[void]$ws.Cells.Replace('Nejapa de Madero', 'Nejapa De Madero')
[void]$ws.Cells.Replace('Oaxaca de Juárez', 'Oaxaca De Juárez')
[void]$ws.Cells.Replace('Ocotlán de Morelos', 'Ocotlán De Morelos')
[void]$ws.Cells.Replace('Pinotepa de Don Luis', 'Pinotepa De Don Luis')
[void]$ws.Cells.Replace('Putla Villa de Guerrero', 'Putla Villa De Guerrero')
[void]$ws.Cells.Replace('Reforma de Pineda', 'Reforma De Pineda')
[void]$ws.Cells.Replace('Rojas de Cuauhtémoc', 'Rojas De Cuauhtémoc')
[void]$ws.Cells.Replace('San Antonino el Alto', 'San Antonino El Alto')
[void]$ws.Cells.Replace('San José del Progreso', 'San José Del Progreso')
[void]$ws.Cells.Replace('San Juan Bautista Lo de Soto', 'San Juan Bautista Lo De Soto')
[void]$ws.Cells.Replace('San Juan del Estado', 'San Juan Del Estado')
[void]$ws.Cells.Replace('San Miguel del Puerto', 'San Miguel Del Puerto')
[void]$ws.Cells.Replace('San Miguel el Grande', 'San Miguel El Grande')
[void]$ws.Cells.Replace('San Pablo Villa de Mitla', 'San Pablo Villa De Mitla')
[void]$ws.Cells.Replace('San Pedro y San Pablo Ayutla', 'San Pedro Y San Pablo Ayutla')
[void]$ws.Cells.Replace('San Pedro y San Pablo Teposcolula', 'San Pedro Y San Pablo Teposcolula')
[void]$ws.Cells.Replace('San Pedro y San Pablo Tequixtepec', 'San Pedro Y San Pablo Tequixtepec')
[void]$ws.Cells.Replace('Santa Cruz de Bravo', 'Santa Cruz De Bravo')
[void]$ws.Cells.Replace('Santa Cruz Tacache de Mina', 'Santa Cruz Tacache De Mina')
[void]$ws.Cells.Replace('Santa Lucía del Camino', 'Santa Lucía Del Camino')
[void]$ws.Cells.Replace('Santa María Jalapa del Marqués', 'Santa María Jalapa Del Marqués')
[void]$ws.Cells.Replace('Tamazulápam del Espíritu Santo', 'Tamazulápam Del Espíritu Santo')
[void]$ws.Cells.Replace('Tanetze de Zaragoza', 'Tanetze De Zaragoza')
[void]$ws.Cells.Replace('Tataltepec de Valdés', 'Tataltepec De Valdés')
[void]$ws.Cells.Replace('Teococuilco de Marcos Pérez', 'Teococuilco De Marcos Pérez')
[void]$ws.Cells.Replace('Teotitlán del Valle', 'Teotitlán Del Valle')
[void]$ws.Cells.Replace('Tepelmeme Villa de Morelos', 'Tepelmeme Villa De Morelos')
[void]$ws.Cells.Replace('Heroica Villa Tezoatlán de Segura y Luna, Cuna de la Independencia de Oaxaca', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca')
[void]$ws.Cells.Replace('Tlacolula de Matamoros', 'Tlacolula De Matamoros')
[void]$ws.Cells.Replace('Tlalixtac de Cabrera', 'Tlalixtac De Cabrera')
[void]$ws.Cells.Replace('Totontepec Villa de Morelos', 'Totontepec Villa De Morelos')
[void]$ws.Cells.Replace('Villa de Chilapa de Díaz', 'Villa De Chilapa De Díaz')
[void]$ws.Cells.Replace('Villa de Etla', 'Villa De Etla')
[void]$ws.Cells.Replace('Villa de Tututepec', 'Villa De Tututepec')
[void]$ws.Cells.Replace('Villa de Zaachila', 'Villa De Zaachila')
[void]$ws.Cells.Replace('Villa Sola de Vega', 'Villa Sola De Vega')
[void]$ws.Cells.Replace('Villa Talea de Castro', 'Villa Talea De Castro')
[void]$ws.Cells.Replace('Zapotitlán del Río', 'Zapotitlán Del Río')
[void]$ws.Cells.Replace('Zimatlán de Álvarez', 'Zimatlán De Álvarez')
[void]$ws.Cells.Replace('Chalchicomula de Sesma', 'Chalchicomula De Sesma')
[void]$ws.Cells.Replace('Chila de la Sal', 'Chila De La Sal')
[void]$ws.Cells.Replace('Cuayuca de Andrade', 'Cuayuca De Andrade')
[void]$ws.Cells.Replace('Cuetzalan del Progreso', 'Cuetzalan Del Progreso')
[void]$ws.Cells.Replace('Huehuetlán el Chico', 'Huehuetlán El Chico')
[void]$ws.Cells.Replace('Huehuetlán el Grande', 'Huehuetlán El Grande')
[void]$ws.Cells.Replace('Huitzilan de Serdán', 'Huitzilan De Serdán')
[void]$ws.Cells.Replace('Ixcamilpa de Guerrero', 'Ixcamilpa De Guerrero')
[void]$ws.Cells.Replace('Izúcar de Matamoros', 'Izúcar De Matamoros')
[void]$ws.Cells.Replace('Los Reyes de Juárez', 'Los Reyes De Juárez')
[void]$ws.Cells.Replace('Mazapiltepec de Juárez', 'Mazapiltepec De Juárez')
[void]$ws.Cells.Replace('Palmar de Bravo', 'Palmar De Bravo')
[void]$ws.Cells.Replace('San Diego la Mesa Tochimiltzingo', 'San Diego La Mesa Tochimiltzingo')
[void]$ws.Cells.Replace('San Nicolás de los Ranchos', 'San Nicolás De Los Ranchos')
[void]$ws.Cells.Replace('San Salvador el Seco', 'San Salvador El Seco')
[void]$ws.Cells.Replace('San Salvador el Verde', 'San Salvador El Verde')
[void]$ws.Cells.Replace('Tecali de Herrera', 'Tecali De Herrera')
[void]$ws.Cells.Replace('Tepanco de López', 'Tepanco De López')
[void]$ws.Cells.Replace('Tepatlaxco de Hidalgo', 'Tepatlaxco De Hidalgo')
[void]$ws.Cells.Replace('Tepexi de Rodríguez', 'Tepexi De Rodríguez')
[void]$ws.Cells.Replace('Tetela de Ocampo', 'Tetela De Ocampo')
[void]$ws.Cells.Replace('Teteles de Avila Castillo', 'Teteles De Avila Castillo')
[void]$ws.Cells.Replace('Tlacotepec de Benito Juárez', 'Tlacotepec De Benito Juárez')
[void]$ws.Cells.Replace('Totoltepec de Guerrero', 'Totoltepec De Guerrero')
[void]$ws.Cells.Replace('Xayacatlán de Bravo', 'Xayacatlán De Bravo')
[void]$ws.Cells.Replace('Amealco de Bonfil', 'Amealco De Bonfil')
[void]$ws.Cells.Replace('Cadereyta de Montes', 'Cadereyta De Montes')
[void]$ws.Cells.Replace('Jalpan de Serra', 'Jalpan De Serra')
[void]$ws.Cells.Replace('Landa de Matamoros', 'Landa De Matamoros')
[void]$ws.Cells.Replace('Pinal de Amoles', 'Pinal De Amoles')
[void]$ws.Cells.Replace('Axtla de Terrazas', 'Axtla De Terrazas')
[void]$ws.Cells.Replace('San Ciro de Acosta', 'San Ciro De Acosta')
[void]$ws.Cells.Replace('Tanquián de Escobedo', 'Tanquián De Escobedo')
[void]$ws.Cells.Replace('Jalpa de Méndez', 'Jalpa De Méndez')
[void]$ws.Cells.Replace('Soto la Marina', 'Soto La Marina')
[void]$ws.Cells.Replace('Acuamanala de Miguel Hidalgo', 'Acuamanala De Miguel Hidalgo')
[void]$ws.Cells.Replace('Amaxac de Guerrero', 'Amaxac De Guerrero')
[void]$ws.Cells.Replace('Apetatitlán de Antonio Carvajal', 'Apetatitlán De Antonio Carvajal')
[void]$ws.Cells.Replace('Contla de Juan Cuamatzi', 'Contla De Juan Cuamatzi')
[void]$ws.Cells.Replace('Ixtacuixtla de Mariano Matamoros', 'Ixtacuixtla De Mariano Matamoros')
[void]$ws.Cells.Replace('Mazatecochco de José María Morelos', 'Mazatecochco De José María Morelos')
[void]$ws.Cells.Replace('Nanacamilpa de Mariano Arista', 'Nanacamilpa De Mariano Arista')
[void]$ws.Cells.Replace('Papalotla de Xicohténcatl', 'Papalotla De Xicohténcatl')
[void]$ws.Cells.Replace('San Pablo del Monte', 'San Pablo Del Monte')
[void]$ws.Cells.Replace('Tepetitla de Lardizábal', 'Tepetitla De Lardizábal')
[void]$ws.Cells.Replace('Tetla de la Solidaridad', 'Tetla De La Solidaridad')
[void]$ws.Cells.Replace('Ziltlaltépec de Trinidad Sánchez Santos', 'Ziltlaltépec De Trinidad Sánchez Santos')
[void]$ws.Cells.Replace('Amatlán de los Reyes', 'Amatlán De Los Reyes')
[void]$ws.Cells.Replace('Boca del Río', 'Boca Del Río')
[void]$ws.Cells.Replace('Cosamaloapan de Carpio', 'Cosamaloapan De Carpio')
[void]$ws.Cells.Replace('Hueyapan de Ocampo', 'Hueyapan De Ocampo')
[void]$ws.Cells.Replace('Huiloapan de Cuauhtémoc', 'Huiloapan De Cuauhtémoc')
[void]$ws.Cells.Replace('Ignacio de la Llave', 'Ignacio De La Llave')
[void]$ws.Cells.Replace('Ixhuacán de los Reyes', 'Ixhuacán De Los Reyes')
[void]$ws.Cells.Replace('Ixhuatlán de Madero', 'Ixhuatlán De Madero')
[void]$ws.Cells.Replace('Ixhuatlán del Sureste', 'Ixhuatlán Del Sureste')
[void]$ws.Cells.Replace('Juchique de Ferrer', 'Juchique De Ferrer')
[void]$ws.Cells.Replace('Las Vigas de Ramírez', 'Las Vigas De Ramírez')
[void]$ws.Cells.Replace('Martínez de la Torre', 'Martínez De La Torre')
[void]$ws.Cells.Replace('Medellín de Bravo', 'Medellín De Bravo')
[void]$ws.Cells.Replace('Ozuluama de Mascareñas', 'Ozuluama De Mascareñas')
[void]$ws.Cells.Replace('Paso de Ovejas', 'Paso De Ovejas')
[void]$ws.Cells.Replace('Paso del Macho', 'Paso Del Macho')
[void]$ws.Cells.Replace('Sayula de Alemán', 'Sayula De Alemán')
[void]$ws.Cells.Replace('Vega de Alatorre', 'Vega De Alatorre')
[void]$ws.Cells.Replace('Zontecomatlán de López y Fuentes', 'Zontecomatlán De López Y Fuentes')
[void]$ws.Cells.Replace('Zozocolco de Hidalgo', 'Zozocolco De Hidalgo')
[void]$ws.Cells.Replace('Noria de Ángeles', 'Noria De Ángeles')

# 3. Fix a floating point rounding difference on D476
$ws.Range("D476").Value = 0.009369369369369367

# 4. Remove the trailing footnote/metadata rows (1275-1279) that sat below
#    the data's final "Total" row (1273); row 1274 was already blank.
$ws.Rows.Item(1279).Delete()
$ws.Rows.Item(1278).Delete()
$ws.Rows.Item(1277).Delete()
$ws.Rows.Item(1276).Delete()
$ws.Rows.Item(1275).Delete()
